$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "ID Competição" column (B) values from 67 to 267 for data rows 2-9
$ws.Range("B2:B9").Value = 267
